# "Add support for MacOS" — update the sample recipient row data so the
# workbook carries a generic MacOS-friendly test contact instead of the
# previous (Windows-authored) sample values.
#
# Column layout (row 1 = header): A=Name B=Number C=FormalMessage
# D=InformalMessage E=FormalSalutation F=InformalSalutation G=Informal
# H=Send I=Status

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name column: replace the sample contact name on every data row ------
$ws.Range("A2:A6").Value = "Vikrant"

# --- Formal message text (row 2) -----------------------------------------
$ws.Range("C2").Value = "how are you?"

# --- Phone number samples --------------------------------------------------
# B2 / B3: plain numeric entry (same shape as before the edit).
$ws.Range("B2").Value = 9289907889
$ws.Range("B3").Value = 9289907889

# B4 / B5: these must stay TEXT (leading/formatted phone-number strings),
# not get auto-converted to numbers. Typing a numeric-looking string into
# .Value/.Formula directly gets reinterpreted as a number, and forcing text
# via NumberFormat="@" (or a "'" quote-prefix) leaves a permanent
# quotePrefix style on the cell. Instead, write the value as a text formula
# (="...") - which evaluates to a *string* result - then copy/paste-special
# as values only; this keeps the cell a plain shared-string text cell with
# no extra style applied, matching how the sheet already stored B4:B6.
$ws.Range("B4").Formula = '="919289907889"'
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B5").Formula = '="+919289907889"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false

# B6 keeps its original text value ("982501851") - nothing to change there.

# --- Selection / cursor position left where the editor left it -----------
$ws.Range("J5").Select() | Out-Null

# --- Row heights: drop the explicit 15.75pt custom height left over from
# the old (Windows) save so rows fall back to the sheet default, matching
# the look of a freshly (Mac) resaved sheet.
$ws.Rows("1:6").AutoFit()
